$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3; this shifts the existing row 3 (Racing Montevideo match) down to row 4
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the new match data (Palmeiras vs Gremio)
$ws.Range('A3').Value = 'rVIQRBBR'
$ws.Range('B3').NumberFormat = '@'
$ws.Range('B3').Value = '08/11/2024'
$ws.Range('B3').Style = 'Normal'
$ws.Range('C3').Value = '21:30'
$ws.Range('D3').Value = 'BRAZIL - SERIE A BETANO'
$ws.Range('E3').Value = 'Palmeiras'
$ws.Range('F3').Value = 'Gremio'
$ws.Range('G3').Value = 1.4
$ws.Range('H3').Value = 4.75
$ws.Range('I3').Value = 7.5
$ws.Range('J3').Value = 1.91
$ws.Range('K3').Value = 2.5
$ws.Range('L3').Value = 6.5
$ws.Range('M3').Value = 1.03
$ws.Range('N3').Value = 15
$ws.Range('O3').Value = 1.2
$ws.Range('P3').Value = 4.5
$ws.Range('Q3').Value = 1.6
$ws.Range('R3').Value = 2.3
$ws.Range('S3').Value = 1.3
$ws.Range('T3').Value = 3.4
$ws.Range('U3').Value = 1.91
$ws.Range('V3').Value = 1.91
$ws.Range('W3').Value = 8
$ws.Range('X3').Value = 7
$ws.Range('Y3').Value = 8.5
$ws.Range('Z3').Value = 9.5
$ws.Range('AA3').Value = 11
$ws.Range('AB3').Value = 23
$ws.Range('AC3').Value = 15
$ws.Range('AD3').Value = 9.5
$ws.Range('AE3').Value = 19
$ws.Range('AF3').Value = 51
$ws.Range('AG3').Value = 19
$ws.Range('AH3').Value = 41
$ws.Range('AI3').Value = 21
$ws.Range('AJ3').Value = 81
$ws.Range('AK3').Value = 51
$ws.Range('AL3').Value = 41
$ws.Range('AM3').Value = 251
$ws.Range('AN3').Value = 3.4
$ws.Range('AO3').Value = 6.5
$ws.Range('AP3').Value = 17
$ws.Range('AQ3').Value = 19
$ws.Range('AR3').Value = 41
$ws.Range('AS3').Value = 101
$ws.Range('AT3').Value = 3.4
$ws.Range('AU3').Value = 8.5
$ws.Range('AV3').Value = 51
$ws.Range('AW3').Value = 8.5
$ws.Range('AX3').Value = 34
$ws.Range('AY3').Value = 34
$ws.Range('AZ3').Value = 126
$ws.Range('BA3').Value = 126
$ws.Range('BB3').Value = 251
$ws.Range('BC3').Value = 126
$ws.Range('BD3').Value = 151
